$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.954.99"
$ws.Range("E2").Value = "  +1.73%  "

$ws.Range("D3").Value = "3.256.49"

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.30"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.07"
$ws.Range("E6").Value = "  +4.40%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("E9").Value = "  +4.23%  "

$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("E11").Value = "  +2.60%  "

$ws.Range("D12").Value = "3.818.46"
$ws.Range("E12").Value = "  +0.71%  "

$ws.Range("E13").Value = "  +0.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.77"
$ws.Range("E14").Value = "  +3.99%  "

$ws.Range("D15").Value = "67.961.26"
$ws.Range("E15").Value = "  +1.75%  "

$ws.Range("E16").Value = "  +1.89%  "

$ws.Range("D17").Value = "3.251.29"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.54"
$ws.Range("E19").Value = "  +2.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "379.09"
$ws.Range("E20").Value = "  +3.72%  "

$ws.Range("E21").Value = "  +3.17%  "

$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.29"
$ws.Range("E23").Value = "  +2.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.514"
$ws.Range("E24").Value = "  +2.04%  "

$ws.Range("E25").Value = "  +1.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.08"
$ws.Range("E26").Value = "  +2.58%  "

$ws.Range("E27").Value = "  +2.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("E29").Value = "  +1.41%  "

$ws.Range("E30").Value = "  +2.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.86"
$ws.Range("E31").Value = "  +2.35%  "

$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.03"
$ws.Range("E33").Value = "  +4.42%  "

$ws.Range("E34").Value = "  +3.92%  "

$ws.Range("E35").Value = "  +6.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.70"
$ws.Range("E36").Value = "  -3.35%  "

$ws.Range("E37").Value = "  +1.67%  "

$ws.Range("E38").Value = "  +0.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.50"
$ws.Range("E39").Value = "  -0.49%  "

$ws.Range("E40").Value = "  +8.28%  "

$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.61"
$ws.Range("E41").Value = "  +0.72%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.61"
$ws.Range("E42").Value = "  +3.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.56"
$ws.Range("E43").Value = "  +4.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "347.48"
$ws.Range("E44").Value = "  +5.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.13"
$ws.Range("E45").Value = "  +1.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0687"
$ws.Range("E46").Value = "  +2.56%  "

$ws.Range("D47").Value = "2.636.67"
$ws.Range("E47").Value = "  -1.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0287"
$ws.Range("E48").Value = "  +3.39%  "

$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("E50").Value = "  +2.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.17"
$ws.Range("E51").Value = "  +3.05%  "
